# Auto-applied numeric updates to Ifrit_Profits workbook (scheduled runner refresh)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 3041.5
$ws.Range("I98").Value = 3255.4443
$ws.Range("J98").Value = 2399.6667
$ws.Range("K98").Value = 3255.4443
$ws.Range("L98").Value = 2399.6667
$ws.Range("M98").Value = -1757.4443
$ws.Range("N98").Value = -5395.6667
$ws.Range("H122").Value = 3041.5
$ws.Range("I122").Value = 3255.4443
$ws.Range("J122").Value = 2399.6667
$ws.Range("K122").Value = 9766.332900000001
$ws.Range("L122").Value = 7199.000100000001
$ws.Range("M122").Value = -7316.332900000001
$ws.Range("N122").Value = -12099.0001
$ws.Range("H125").Value = 1074.4445
$ws.Range("I125").Value = 683.75
$ws.Range("J125").Value = 4200
$ws.Range("K125").Value = 6153.75
$ws.Range("L125").Value = 37800
$ws.Range("M125").Value = -3693.75
$ws.Range("N125").Value = -42720

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 517.8
$ws.Range("I97").Value = 517.8
$ws.Range("K97").Value = 517.8
$ws.Range("M97").Value = -21.79999999999995
$ws.Range("H109").Value = 37500
$ws.Range("J109").Value = 37500
$ws.Range("L109").Value = 37500
$ws.Range("N109").Value = -40274
$ws.Range("H110").Value = 733.3333
$ws.Range("I110").Value = 900
$ws.Range("J110").Value = 400
$ws.Range("K110").Value = 900
$ws.Range("L110").Value = 400
$ws.Range("M110").Value = 1145
$ws.Range("N110").Value = -4490
$ws.Range("H112").Value = 31064.334
$ws.Range("J112").Value = 31064.334
$ws.Range("L112").Value = 31064.334
$ws.Range("N112").Value = -34018.334
$ws.Range("H124").Value = 8607.25
$ws.Range("J124").Value = 8607.25
$ws.Range("L124").Value = 8607.25
$ws.Range("N124").Value = -18427.25
$ws.Range("H125").Value = 46392
$ws.Range("J125").Value = 46392
$ws.Range("L125").Value = 46392
$ws.Range("N125").Value = -56232
$ws.Range("H132").Value = 2477.25
$ws.Range("I132").Value = 2624.05
$ws.Range("J132").Value = 2293.75
$ws.Range("K132").Value = 7872.150000000001
$ws.Range("L132").Value = 6881.25
$ws.Range("M132").Value = -5342.150000000001
$ws.Range("N132").Value = -11941.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1999.75
$ws.Range("I94").Value = 1999.75
$ws.Range("K94").Value = 1999.75
$ws.Range("M94").Value = -1548.75
$ws.Range("H107").Value = 2154.9092
$ws.Range("I107").Value = 2086.375
$ws.Range("K107").Value = 2086.375
$ws.Range("M107").Value = -166.375
$ws.Range("H134").Value = 53662.74
$ws.Range("I134").Value = 68024.72
$ws.Range("J134").Value = 1959.6
$ws.Range("K134").Value = 204074.16
$ws.Range("L134").Value = 5878.799999999999
$ws.Range("M134").Value = -201539.16
$ws.Range("N134").Value = -10948.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2487.2727
$ws.Range("I31").Value = 1726.6666
$ws.Range("J31").Value = 3400
$ws.Range("K31").Value = 1726.6666
$ws.Range("L31").Value = 3400
$ws.Range("M31").Value = -1431.6666
$ws.Range("N31").Value = -3990
$ws.Range("H34").Value = 2487.2727
$ws.Range("I34").Value = 1726.6666
$ws.Range("J34").Value = 3400
$ws.Range("K34").Value = 1726.6666
$ws.Range("L34").Value = 3400
$ws.Range("M34").Value = -1524.6666
$ws.Range("N34").Value = -3804
$ws.Range("H99").Value = 1158.64
$ws.Range("I99").Value = 1117.7858
$ws.Range("J99").Value = 1210.6364
$ws.Range("K99").Value = 1117.7858
$ws.Range("L99").Value = 1210.6364
$ws.Range("M99").Value = 380.2141999999999
$ws.Range("N99").Value = -4206.6364
$ws.Range("H126").Value = 1158.64
$ws.Range("I126").Value = 1117.7858
$ws.Range("J126").Value = 1210.6364
$ws.Range("K126").Value = 3353.3574
$ws.Range("L126").Value = 3631.9092
$ws.Range("M126").Value = -883.3574000000003
$ws.Range("N126").Value = -8571.9092

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 1575
$ws.Range("J74").Value = 3300
$ws.Range("L74").Value = 9900
$ws.Range("N74").Value = -12022
$ws.Range("H77").Value = 1575
$ws.Range("J77").Value = 3300
$ws.Range("L77").Value = 29700
$ws.Range("N77").Value = -40308
$ws.Range("H96").Value = 35354950
$ws.Range("J96").Value = 35354950
$ws.Range("L96").Value = 106064850
$ws.Range("N96").Value = -106068968
$ws.Range("H105").Value = 181602400
$ws.Range("J105").Value = 181602400
$ws.Range("L105").Value = 544807200
$ws.Range("N105").Value = -544812442
$ws.Range("H110").Value = 3062.5
$ws.Range("J110").Value = 3740
$ws.Range("L110").Value = 11220
$ws.Range("N110").Value = -19400
$ws.Range("H133").Value = 5023.839
$ws.Range("J133").Value = 7700
$ws.Range("L133").Value = 23100
$ws.Range("N133").Value = -33220

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1504.3334
$ws.Range("I102").Value = 1017.7143
$ws.Range("J102").Value = 2185.6
$ws.Range("K102").Value = 1017.7143
$ws.Range("L102").Value = 2185.6
$ws.Range("M102").Value = 604.2857
$ws.Range("N102").Value = -5429.6
$ws.Range("H122").Value = 5977.0557
$ws.Range("I122").Value = 6255.4375
$ws.Range("J122").Value = 3750
$ws.Range("K122").Value = 18766.3125
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -16316.3125
$ws.Range("N122").Value = -16150
$ws.Range("H126").Value = 2372.3635
$ws.Range("I126").Value = 3003.6667
$ws.Range("J126").Value = 1614.8
$ws.Range("K126").Value = 9011.000100000001
$ws.Range("L126").Value = 4844.4
$ws.Range("M126").Value = -6541.000100000001
$ws.Range("N126").Value = -9784.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2418.4375
$ws.Range("I7").Value = 2154.182
$ws.Range("J7").Value = 2999.8
$ws.Range("K7").Value = 2154.182
$ws.Range("L7").Value = 2999.8
$ws.Range("M7").Value = -2042.182
$ws.Range("N7").Value = -3223.8
$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 800
$ws.Range("K22").Value = 800
$ws.Range("M22").Value = -505
$ws.Range("H27").Value = 600
$ws.Range("I27").Value = 800
$ws.Range("K27").Value = 800
$ws.Range("M27").Value = -693
$ws.Range("H40").Value = 1756.9688
$ws.Range("I40").Value = 1725.7916
$ws.Range("J40").Value = 1850.5
$ws.Range("K40").Value = 1725.7916
$ws.Range("L40").Value = 1850.5
$ws.Range("M40").Value = -1589.7916
$ws.Range("N40").Value = -2122.5
$ws.Range("H93").Value = 2220.077
$ws.Range("I93").Value = 1678.6842
$ws.Range("K93").Value = 1678.6842
$ws.Range("M93").Value = -430.6841999999999
$ws.Range("H122").Value = 8585.611000000001
$ws.Range("I122").Value = 11020
$ws.Range("J122").Value = 3716.8333
$ws.Range("K122").Value = 33060
$ws.Range("L122").Value = 11150.4999
$ws.Range("M122").Value = -30610
$ws.Range("N122").Value = -16050.4999
$ws.Range("H126").Value = 2418.4375
$ws.Range("I126").Value = 2154.182
$ws.Range("J126").Value = 2999.8
$ws.Range("K126").Value = 6462.545999999999
$ws.Range("L126").Value = 8999.400000000001
$ws.Range("M126").Value = -3992.545999999999
$ws.Range("N126").Value = -13939.4
$ws.Range("H127").Value = 47847.5
$ws.Range("J127").Value = 47847.5
$ws.Range("L127").Value = 47847.5
$ws.Range("N127").Value = -57767.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 30724.75
$ws.Range("J104").Value = 30724.75
$ws.Range("L104").Value = 30724.75
$ws.Range("N104").Value = -37712.75
$ws.Range("H122").Value = 2130.625
$ws.Range("I122").Value = 1016.6667
$ws.Range("J122").Value = 2799
$ws.Range("K122").Value = 3050.0001
$ws.Range("L122").Value = 8397
$ws.Range("M122").Value = -600.0001000000002
$ws.Range("N122").Value = -13297
$ws.Range("H126").Value = 2422.524
$ws.Range("I126").Value = 2514.3684
$ws.Range("J126").Value = 1550
$ws.Range("K126").Value = 7543.1052
$ws.Range("L126").Value = 4650
$ws.Range("M126").Value = -5073.1052
$ws.Range("N126").Value = -9590
$ws.Range("H136").Value = 1317.7391
$ws.Range("I136").Value = 1392.9231
$ws.Range("J136").Value = 1220
$ws.Range("K136").Value = 4178.7693
$ws.Range("L136").Value = 3660
$ws.Range("M136").Value = -1628.7693
$ws.Range("N136").Value = -8760
